# Expenditures.xlsx update: add the Digi-Key (board components) order row
# to the "Mike" order table (columns M:O), row 6.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the date first, then the item description, so the shared-string
# table grows in the same order as the authored workbook (Oct.31 before
# the Board Components description).
$ws.Range("O6").Value = "Oct.31"
$ws.Range("M6").Value = "Board Components (See Component Invoice.xlsx)"
$ws.Range("N6").Value = 64.93

# Match the currency formatting used by the other order-cost cells in this
# column, and wrap the long item description / cost text like the rest of
# the order rows that carry multi-line content.
$ws.Range("N6").NumberFormat = $ws.Range("N4").NumberFormat
$ws.Range("M6").WrapText = $true
$ws.Range("N6").WrapText = $true

# The extra-tall row needed to show the wrapped description.
$ws.Rows.Item(6).RowHeight = 44.25

# Leave the selection where the author left it after entering this data.
$ws.Range("M6").Select() | Out-Null
